$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7820801138877869
$ws.Range("B1").Value = 1.59805166721344
$ws.Range("C1").Value = 1.67356538772583
$ws.Range("D1").Value = 6.398558616638184
$ws.Range("E1").Value = 3.561268091201782
